$wb = $excel.ActiveWorkbook

# Rename the three tabs to their new destinations
$wb.Worksheets.Item(1).Name = "Orden por ingresos"
$wb.Worksheets.Item(2).Name = "Orden por Ingreso siendo hombre"
$wb.Worksheets.Item(3).Name = "Orden por ingreso siendo mujer"

# Update the remembered selection on "Orden por ingreso siendo mujer" (was F6 -> F9)
$wb.Worksheets.Item(3).Range("F9").Select()

# Update the remembered selection on "Orden por ingresos" (was C6 -> C4);
# select this sheet last so it stays the active/tabSelected sheet
$wb.Worksheets.Item(1).Range("C4").Select()
